$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the NIVALDO row (account 005277762, balance 64450) entirely.
$ws.Rows.Item(8).Delete()

# Remove the ALBERTO row (account 004480970, balance 16352.97) entirely.
# After the previous deletion it has shifted up from row 14 to row 13.
$ws.Rows.Item(13).Delete()

# The WALQUIRIA row (account 005103059) still sits at row 4 with its old
# balance. Remove it from there; it will be re-inserted further down with
# its updated balance.
$ws.Rows.Item(4).Delete()

# Re-insert the WALQUIRIA row right after JOAQUIM (and before THIAGO) with
# its new balance of 30937.03.
$ws.Rows.Item(10).Insert()
$ws.Cells.Item(10, 1).NumberFormat = "@"
$ws.Cells.Item(10, 1).Value = "005103059"
$ws.Cells.Item(10, 2).Value = "WALQUIRIA"
$ws.Cells.Item(10, 3).Value = 30937.03
